# Generate Report for Handoff
# Updates status/datetime/error-detail fields across the Overview, zh-cn and
# de-de sheets, plus a few column-width tweaks.

$wb = $excel.ActiveWorkbook

$msg2df = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fc28ae4ee5a9d616e3fb8a4280c7f7b26e0bac82/e2e/2df024e5-f384-4970-a0a5-31d3bfdf89a3.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f4c6e4123a66c24a3a5bd2920b3eccd0fa006d82/e2e/2df024e5-f384-4970-a0a5-31d3bfdf89a3.md."
$msgE496 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fc28ae4ee5a9d616e3fb8a4280c7f7b26e0bac82/e2e/e49630ad-40ac-44e9-a236-cf67916150cf.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f4c6e4123a66c24a3a5bd2920b3eccd0fa006d82/e2e/e49630ad-40ac-44e9-a236-cf67916150cf.md."

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("E2").Value = "In Translation"
$ws.Range("F2").Value = "In Translation"
$ws.Range("G2").Value = "2016-10-21 01:12:00"

$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-10-21 01:12:00"

$ws.Columns.Item(5).ColumnWidth = 16.333333333333336
$ws.Columns.Item(6).ColumnWidth = 16.333333333333336

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("C2").Value = "In Translation"
$ws.Range("H2").Value = "2016-10-21 01:11:48"
$ws.Range("P2").Value = $msg2df

$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("H3").Value = "2016-10-21 01:11:48"
$ws.Range("P3").Value = $msgE496

$ws.Columns.Item(3).ColumnWidth = 16.333333333333336
$ws.Columns.Item(16).ColumnWidth = 39.16666666666667

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("C2").Value = "In Translation"
$ws.Range("H2").Value = "2016-10-21 01:12:00"
$ws.Range("P2").Value = $msg2df

$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("H3").Value = "2016-10-21 01:12:00"
$ws.Range("P3").Value = $msgE496

$ws.Columns.Item(3).ColumnWidth = 16.333333333333336
$ws.Columns.Item(16).ColumnWidth = 39.16666666666667
